$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be auto-parsed as numbers
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "34.123.29"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "1.790.73"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +3.05%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").Value = "30.38"
$ws.Range("E8").Value = "  -3.80%  "
$ws.Range("D9").Value = "46.51"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").Value = "0.0668"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "0.0924"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").Value = "2.048.92"
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").Value = "1.814.42"
$ws.Range("E14").Value = "  +3.34%  "
$ws.Range("D15").Value = "0.627"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "34.135.03"
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("D17").Value = "10.31"
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("D18").Value = "4.20"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").Value = "68.96"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").Value = "252.52"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "10.36"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "4.22"
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D25").Value = "2.12"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("D26").Value = "158.54"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "16.53"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").Value = "3.84"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("E34").Value = "  +4.14%  "
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("D36").Value = "1.504.69"
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("D37").Value = "1.07"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").Value = "83.63"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "2.35"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("D43").Value = "0.905"
$ws.Range("E43").Value = "  +4.17%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.0516"
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "2.05"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "1.947.89"
$ws.Range("E47").Value = "  +2.49%  "
$ws.Range("D48").Value = "5.74"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("D50").Value = "11.82"
$ws.Range("E50").Value = "  +6.58%  "
$ws.Range("D51").Value = "51.51"

Write-Output "Applied cryptos update"